$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1887.8889
$ws.Range("I32").Value = 1917.8
$ws.Range("J32").Value = 1850.5
$ws.Range("K32").Value = 1917.8
$ws.Range("L32").Value = 1850.5
$ws.Range("M32").Value = -1591.8
$ws.Range("N32").Value = -2502.5

$ws.Range("H51").Value = 12799.667
$ws.Range("I51").Value = 9449
$ws.Range("J51").Value = 14475
$ws.Range("K51").Value = 9449
$ws.Range("L51").Value = 14475
$ws.Range("M51").Value = -8965
$ws.Range("N51").Value = -15443

$ws.Range("H75").Value = 56000
$ws.Range("J75").Value = 56000
$ws.Range("L75").Value = 56000
$ws.Range("N75").Value = -57872

$ws.Range("H78").Value = 56000
$ws.Range("J78").Value = 56000
$ws.Range("L78").Value = 168000
$ws.Range("N78").Value = -177360

$ws.Range("H100").Value = 2189.6
$ws.Range("I100").Value = 2236.4211
$ws.Range("J100").Value = 1300
$ws.Range("K100").Value = 2236.4211
$ws.Range("L100").Value = 1300
$ws.Range("M100").Value = -1695.4211
$ws.Range("N100").Value = -2382

$ws.Range("H107").Value = 344.22223
$ws.Range("I107").Value = 286.06668
$ws.Range("J107").Value = 635
$ws.Range("K107").Value = 286.06668
$ws.Range("L107").Value = 635
$ws.Range("M107").Value = 1633.93332
$ws.Range("N107").Value = -4475

$ws.Range("H113").Value = 3499.75
$ws.Range("I113").Value = 3499.6667
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3499.6667
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -245.6667000000002
$ws.Range("N113").Value = -10008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 71166.664
$ws.Range("J56").Value = 71166.664
$ws.Range("L56").Value = 71166.664
$ws.Range("N56").Value = -72650.664

$ws.Range("H132").Value = 4178.8
$ws.Range("I132").Value = 4073.5
$ws.Range("J132").Value = 4600
$ws.Range("K132").Value = 12220.5
$ws.Range("L132").Value = 13800
$ws.Range("M132").Value = -9690.5
$ws.Range("N132").Value = -18860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1170.25
$ws.Range("I16").Value = 1090.5
$ws.Range("K16").Value = 1090.5
$ws.Range("M16").Value = -803.5

$ws.Range("H58").Value = 8749.818
$ws.Range("I58").Value = 8319.857
$ws.Range("J58").Value = 9502.25
$ws.Range("K58").Value = 8319.857
$ws.Range("L58").Value = 9502.25
$ws.Range("M58").Value = -8116.857
$ws.Range("N58").Value = -9908.25

$ws.Range("H94").Value = 2869.5
$ws.Range("I94").Value = 2650.3
$ws.Range("K94").Value = 2650.3
$ws.Range("M94").Value = -2199.3

$ws.Range("H105").Value = 2360
$ws.Range("I105").Value = 1800
$ws.Range("K105").Value = 1800
$ws.Range("M105").Value = -53

$ws.Range("H113").Value = 1170.25
$ws.Range("I113").Value = 1090.5
$ws.Range("K113").Value = 1090.5
$ws.Range("M113").Value = 1079.5

$ws.Range("H134").Value = 4099.364
$ws.Range("I134").Value = 3843.6667
$ws.Range("J134").Value = 5250
$ws.Range("K134").Value = 11531.0001
$ws.Range("L134").Value = 15750
$ws.Range("M134").Value = -8996.000100000001
$ws.Range("N134").Value = -20820

$ws.Range("H136").Value = 8749.818
$ws.Range("I136").Value = 8319.857
$ws.Range("J136").Value = 9502.25
$ws.Range("K136").Value = 24959.571
$ws.Range("L136").Value = 28506.75
$ws.Range("M136").Value = -22409.571
$ws.Range("N136").Value = -33606.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1946.7142
$ws.Range("J12").Value = 2711.8
$ws.Range("L12").Value = 8135.400000000001
$ws.Range("N12").Value = -8481.400000000001

$ws.Range("H61").Value = 270.83334
$ws.Range("I61").Value = 93.75
$ws.Range("K61").Value = 281.25
$ws.Range("M61").Value = -66.25

$ws.Range("H80").Value = 7500
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 45000
$ws.Range("N80").Value = -46872

$ws.Range("H83").Value = 7500
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 135000
$ws.Range("N83").Value = -144360

$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992

$ws.Range("H134").Value = 5099
$ws.Range("I134").Value = 5099
$ws.Range("K134").Value = 15297
$ws.Range("M134").Value = -10227

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 29596
$ws.Range("J136").Value = 29596
$ws.Range("L136").Value = 88788
$ws.Range("N136").Value = -93888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1876.5834
$ws.Range("I22").Value = 941.8333
$ws.Range("J22").Value = 2811.3333
$ws.Range("K22").Value = 941.8333
$ws.Range("L22").Value = 2811.3333
$ws.Range("M22").Value = -646.8333
$ws.Range("N22").Value = -3401.3333

$ws.Range("H27").Value = 1876.5834
$ws.Range("I27").Value = 941.8333
$ws.Range("J27").Value = 2811.3333
$ws.Range("K27").Value = 941.8333
$ws.Range("L27").Value = 2811.3333
$ws.Range("M27").Value = -834.8333
$ws.Range("N27").Value = -3025.3333

$ws.Range("H46").Value = 4406.1875
$ws.Range("I46").Value = 4366.6665
$ws.Range("J46").Value = 4429.9
$ws.Range("K46").Value = 4366.6665
$ws.Range("L46").Value = 4429.9
$ws.Range("M46").Value = -4178.6665
$ws.Range("N46").Value = -4805.9

$ws.Range("H55").Value = 840.4375
$ws.Range("I55").Value = 422.5
$ws.Range("K55").Value = 422.5
$ws.Range("M55").Value = -249.5

$ws.Range("H132").Value = 3667.1765
$ws.Range("I132").Value = 3549.3845
$ws.Range("K132").Value = 10648.1535
$ws.Range("M132").Value = -8118.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 714.86664
$ws.Range("I107").Value = 374.9091
$ws.Range("J107").Value = 1649.75
$ws.Range("K107").Value = 1124.7273
$ws.Range("L107").Value = 4949.25
$ws.Range("M107").Value = 795.2727
$ws.Range("N107").Value = -8789.25

$ws.Range("H113").Value = 800
$ws.Range("J113").Value = 900
$ws.Range("L113").Value = 2700
$ws.Range("N113").Value = -7040

$ws.Range("H136").Value = 3515.75
$ws.Range("I136").Value = 3229.5
$ws.Range("K136").Value = 9688.5
$ws.Range("M136").Value = -7138.5
